# Update stock data: a new trade was recorded, so the existing trade-log
# rows shift down by one and a new row is written at the top of the data
# (row 5, directly below the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# --- Shift the existing rows down by one (write literal values rather than
# copying from the source cell, since the DATE column uses a date-formatted
# style and round-tripping its .Value through another cell is unreliable) ---

# Row 8 <- old row 7 (300 @ 37.25, 2025-12-26)
$ws.Cells.Item(8, 1).Value = 46017
$ws.Cells.Item(8, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(8, 2).Value = "NSE"
$ws.Cells.Item(8, 3).Value = "Buy"
$ws.Cells.Item(8, 4).Value = 300
$ws.Cells.Item(8, 5).Value = 37.25
$ws.Cells.Item(8, 6).Value = 11175
$ws.Cells.Item(8, 7).Value = "~"
$ws.Range("J8").Formula = "=Index!`$C`$2"

# Row 7 <- old row 6 (100 @ 36.01, 2026-01-27)
$ws.Cells.Item(7, 1).Value = 46049
$ws.Cells.Item(7, 2).Value = "NSE"
$ws.Cells.Item(7, 3).Value = "Buy"
$ws.Cells.Item(7, 4).Value = 100
$ws.Cells.Item(7, 5).Value = 36.01
$ws.Cells.Item(7, 6).Value = 3601
$ws.Cells.Item(7, 7).Value = "~"
$ws.Range("J7").Formula = "=Index!`$C`$2"

# Row 6 <- old row 5 (100 @ 37.1, 2026-02-06)
$ws.Cells.Item(6, 1).Value = 46059
$ws.Cells.Item(6, 2).Value = "NSE"
$ws.Cells.Item(6, 3).Value = "Buy"
$ws.Cells.Item(6, 4).Value = 100
$ws.Cells.Item(6, 5).Value = 37.1
$ws.Cells.Item(6, 6).Value = 3710
$ws.Cells.Item(6, 7).Value = "~"
$ws.Range("J6").Formula = "=Index!`$C`$2"

# --- New row 5: latest trade ---
$ws.Cells.Item(5, 1).Value = 46062
$ws.Cells.Item(5, 2).Value = "NSE"
$ws.Cells.Item(5, 3).Value = "Buy"
$ws.Cells.Item(5, 4).Value = 100
$ws.Cells.Item(5, 5).Value = 37.8099
$ws.Cells.Item(5, 6).Value = 3807.89
$ws.Cells.Item(5, 7).Value = "CN#252611665409"
$ws.Cells.Item(5, 8).Value = 3.78
$ws.Cells.Item(5, 9).Value = 23.12
$ws.Range("J5").Formula = "=Index!`$C`$2"
